$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3806.842
$ws.Range("I51").Value = 3424.8572
$ws.Range("J51").Value = 4029.6667
$ws.Range("K51").Value = 3424.8572
$ws.Range("L51").Value = 4029.6667
$ws.Range("M51").Value = -2940.8572
$ws.Range("N51").Value = -4997.6667
$ws.Range("H58").Value = 1476.25
$ws.Range("J58").Value = 3999.75
$ws.Range("L58").Value = 11999.25
$ws.Range("N58").Value = -12299.25
$ws.Range("H100").Value = 1944.8182
$ws.Range("I100").Value = 1173.625
$ws.Range("J100").Value = 4001.3333
$ws.Range("K100").Value = 1173.625
$ws.Range("L100").Value = 4001.3333
$ws.Range("M100").Value = -632.625
$ws.Range("N100").Value = -5083.3333
$ws.Range("H103").Value = 545.5
$ws.Range("J103").Value = 545.5
$ws.Range("L103").Value = 1636.5
$ws.Range("N103").Value = -2808.5
$ws.Range("H135").Value = 7814055
$ws.Range("I135").Value = 1009.7727
$ws.Range("J135").Value = 25002754
$ws.Range("K135").Value = 9087.954299999999
$ws.Range("L135").Value = 225024786
$ws.Range("M135").Value = -6552.954299999999
$ws.Range("N135").Value = -225029856
$ws.Range("H138").Value = 7111.0786
$ws.Range("J138").Value = 8898.694
$ws.Range("L138").Value = 26696.082
$ws.Range("N138").Value = -36976.08199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 71434296
$ws.Range("I61").Value = 111114184
$ws.Range("J61").Value = 10500
$ws.Range("K61").Value = 111114184
$ws.Range("L61").Value = 10500
$ws.Range("M61").Value = -111113972
$ws.Range("N61").Value = -10924
$ws.Range("H74").Value = 143018510
$ws.Range("I74").Value = 143018510
$ws.Range("K74").Value = 143018510
$ws.Range("M74").Value = -143017636
$ws.Range("H77").Value = 143018510
$ws.Range("I77").Value = 143018510
$ws.Range("K77").Value = 715092550
$ws.Range("M77").Value = -715088182
$ws.Range("H97").Value = 1685.7
$ws.Range("I97").Value = 720.8570999999999
$ws.Range("K97").Value = 720.8570999999999
$ws.Range("M97").Value = -224.8570999999999
$ws.Range("H132").Value = 50003730
$ws.Range("I132").Value = 3455.2942
$ws.Range("J132").Value = 333338600
$ws.Range("K132").Value = 10365.8826
$ws.Range("L132").Value = 1000015800
$ws.Range("M132").Value = -7835.882599999999
$ws.Range("N132").Value = -1000020860
$ws.Range("H135").Value = 41429
$ws.Range("J135").Value = 41429
$ws.Range("L135").Value = 41429
$ws.Range("N135").Value = -51569
$ws.Range("H136").Value = 71434296
$ws.Range("I136").Value = 111114184
$ws.Range("J136").Value = 10500
$ws.Range("K136").Value = 333342552
$ws.Range("L136").Value = 31500
$ws.Range("M136").Value = -333340002
$ws.Range("N136").Value = -36600
$ws.Range("H137").Value = 29996.5
$ws.Range("J137").Value = 29996.5
$ws.Range("L137").Value = 29996.5
$ws.Range("N137").Value = -40196.5
$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 2719.1428
$ws.Range("I36").Value = 2722.3333
$ws.Range("J36").Value = 2700
$ws.Range("K36").Value = 2722.3333
$ws.Range("L36").Value = 2700
$ws.Range("M36").Value = -2188.3333
$ws.Range("N36").Value = -3768
$ws.Range("H88").Value = 35682.668
$ws.Range("J88").Value = 35682.668
$ws.Range("L88").Value = 35682.668
$ws.Range("N88").Value = -36494.668
$ws.Range("H91").Value = 35682.668
$ws.Range("J91").Value = 35682.668
$ws.Range("L91").Value = 35682.668
$ws.Range("N91").Value = -38490.668
$ws.Range("H134").Value = 1478.2181
$ws.Range("I134").Value = 1478.2181
$ws.Range("K134").Value = 4434.6543
$ws.Range("M134").Value = -1899.6543

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 53849.855
$ws.Range("I60").Value = 14499.5
$ws.Range("J60").Value = 69590
$ws.Range("K60").Value = 14499.5
$ws.Range("L60").Value = 69590
$ws.Range("M60").Value = -13988.5
$ws.Range("N60").Value = -70612

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1574.3684
$ws.Range("J5").Value = 2853.7273
$ws.Range("L5").Value = 8561.1819
$ws.Range("N5").Value = -8785.1819
$ws.Range("H11").Value = 154100.69
$ws.Range("I11").Value = 200082
$ws.Range("J11").Value = 829.6667
$ws.Range("K11").Value = 600246
$ws.Range("L11").Value = 2489.0001
$ws.Range("M11").Value = -600106
$ws.Range("N11").Value = -2769.0001
$ws.Range("H39").Value = 4238.778
$ws.Range("I39").Value = 3125
$ws.Range("J39").Value = 5129.8
$ws.Range("K39").Value = 9375
$ws.Range("L39").Value = 15389.4
$ws.Range("M39").Value = -9081
$ws.Range("N39").Value = -15977.4
$ws.Range("H59").Value = 4831.3335
$ws.Range("J59").Value = 4750
$ws.Range("L59").Value = 14250
$ws.Range("N59").Value = -15330
$ws.Range("H98").Value = 3538.5
$ws.Range("J98").Value = 3962.5557
$ws.Range("L98").Value = 11887.6671
$ws.Range("N98").Value = -14883.6671
$ws.Range("H107").Value = 729
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H129").Value = 4666.6333
$ws.Range("I129").Value = 10840
$ws.Range("J129").Value = 3980.7036
$ws.Range("K129").Value = 32520
$ws.Range("L129").Value = 11942.1108
$ws.Range("M129").Value = -27520
$ws.Range("N129").Value = -21942.1108
$ws.Range("H132").Value = 5559330
$ws.Range("J132").Value = 6064360
$ws.Range("L132").Value = 54579240
$ws.Range("N132").Value = -54584300
$ws.Range("H135").Value = 1574.3684
$ws.Range("J135").Value = 2853.7273
$ws.Range("L135").Value = 25683.5457
$ws.Range("N135").Value = -30753.5457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2767
$ws.Range("I122").Value = 2574
$ws.Range("K122").Value = 7722
$ws.Range("M122").Value = -5272
$ws.Range("H138").Value = 45333.332
$ws.Range("J138").Value = 40000
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4581.325
$ws.Range("I7").Value = 4228.467
$ws.Range("J7").Value = 4793.04
$ws.Range("K7").Value = 4228.467
$ws.Range("L7").Value = 4793.04
$ws.Range("M7").Value = -4116.467
$ws.Range("N7").Value = -5017.04
$ws.Range("H22").Value = 2877.923
$ws.Range("J22").Value = 3351.4
$ws.Range("L22").Value = 3351.4
$ws.Range("N22").Value = -3941.4
$ws.Range("H27").Value = 2877.923
$ws.Range("J27").Value = 3351.4
$ws.Range("L27").Value = 3351.4
$ws.Range("N27").Value = -3565.4
$ws.Range("H46").Value = 1481.5264
$ws.Range("J46").Value = 3662.25
$ws.Range("L46").Value = 3662.25
$ws.Range("N46").Value = -4038.25
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H55").Value = 1100.4166
$ws.Range("I55").Value = 745
$ws.Range("J55").Value = 2166.6667
$ws.Range("K55").Value = 745
$ws.Range("L55").Value = 2166.6667
$ws.Range("M55").Value = -572
$ws.Range("N55").Value = -2512.6667
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()
$ws.Range("H100").Value = 3416.3667
$ws.Range("J100").Value = 3646.2856
$ws.Range("L100").Value = 3646.2856
$ws.Range("N100").Value = -4728.2856
$ws.Range("H122").Value = 5591.1816
$ws.Range("I122").Value = 4967.6665
$ws.Range("K122").Value = 14902.9995
$ws.Range("M122").Value = -12452.9995
$ws.Range("H126").Value = 4581.325
$ws.Range("I126").Value = 4228.467
$ws.Range("J126").Value = 4793.04
$ws.Range("K126").Value = 12685.401
$ws.Range("L126").Value = 14379.12
$ws.Range("M126").Value = -10215.401
$ws.Range("N126").Value = -19319.12

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("I75").Value = 69000
$ws.Range("J75").Value = 79998
$ws.Range("K75").Value = 69000
$ws.Range("L75").Value = 79998
$ws.Range("M75").Value = -68064
$ws.Range("N75").Value = -81870
$ws.Range("I78").Value = 69000
$ws.Range("J78").Value = 79998
$ws.Range("K78").Value = 207000
$ws.Range("L78").Value = 239994
$ws.Range("M78").Value = -202320
$ws.Range("N78").Value = -249354
$ws.Range("H96").Value = 8103.8
$ws.Range("I96").Value = 7849.6665
$ws.Range("J96").Value = 8212.714
$ws.Range("K96").Value = 7849.6665
$ws.Range("L96").Value = 8212.714
$ws.Range("M96").Value = -6476.6665
$ws.Range("N96").Value = -10958.714
$ws.Range("H112").Value = 25387
$ws.Range("J112").Value = 25387
$ws.Range("L112").Value = 25387
$ws.Range("N112").Value = -28341
$ws.Range("H126").Value = 7291
$ws.Range("J126").Value = 4799.8335
$ws.Range("L126").Value = 14399.5005
$ws.Range("N126").Value = -19339.5005
